# RHODE_ISLAND_2017.xlsx cleanup script
# - Rename header row to short English-ish field names
# - Capitalize "de"/"del"/"el"/"los"/"la" particles in municipality /
#   state names (title-case the connector words)
# - Remove stray trailing metadata rows (sample size / source footers)
#   that were left over from a prior sheet beyond the actual data table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case connector words in place names ---
$ws.Range("B3").Value = "Comitán De Domínguez"
$ws.Range("A12").Value = "Ciudad De México"
$ws.Range("A20").Value = "Estado De México"
$ws.Range("B20").Value = "Naucalpan De Juárez"
$ws.Range("B24").Value = "Apaseo El Alto"
$ws.Range("B29").Value = "Acapulco De Juárez"
$ws.Range("B30").Value = "Atoyac De Álvarez"
$ws.Range("B31").Value = "Ayutla De Los Libres"
$ws.Range("B33").Value = "Chilapa De Álvarez"
$ws.Range("B35").Value = "Coyuca De Catalán"
$ws.Range("B36").Value = "Iguala De La Independencia"
$ws.Range("B37").Value = "Zihuatanejo De Azueta"
$ws.Range("B42").Value = "Técpan De Galeana"
$ws.Range("B43").Value = "Tlapa De Comonfort"
$ws.Range("B49").Value = "Huasca De Ocampo"
$ws.Range("B52").Value = "Jacala De Ledezma"
$ws.Range("B59").Value = "Ixtlahuacán Del Río"
$ws.Range("B60").Value = "Lagos De Moreno"
$ws.Range("B62").Value = "Tlajomulco De Zúñiga"
$ws.Range("B81").Value = "San Dionisio Del Mar"
$ws.Range("B93").Value = "Izúcar De Matamoros"
$ws.Range("B96").Value = "San Nicolás De Los Ranchos"
$ws.Range("B105").Value = "Jalpan De Serra"
$ws.Range("B106").Value = "Landa De Matamoros"
$ws.Range("B107").Value = "Pinal De Amoles"
$ws.Range("B111").Value = "Villa De Ramos"
$ws.Range("B132").Value = "Poza Rica De Hidalgo"

# --- Drop the leftover footer rows past the real table (144:148 and
#     476:480), without shifting remaining rows up ---
$ws.Range("A144:D148").ClearContents()
$ws.Range("A476:D480").ClearContents()
